# İş Takip Güncellemesi - 05.08.2025 14:45:13
# Update the "Güncelleme" sheet: move several worked dates from one
# milestone column to another and clear a few "Yapıldı" flags whose
# related work reverted to pending.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Güncelleme")

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $cell = $ws.Range($CellRef)
    # Force a text number format before assigning, otherwise Excel
    # auto-converts a date-shaped string like "2024-11-05" into a date
    # serial number. ClearFormats() afterwards drops the temporary text
    # format again so the cell keeps its original (default, unstyled)
    # look, matching the rest of the sheet - while keeping the cell
    # itself present (even when blanked out) instead of being dropped.
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.ClearFormats()
}

# Row 2
Set-TextValue "P2" ""

# Row 3
Set-TextValue "I3" "2024-11-05"
Set-TextValue "J3" ""
Set-TextValue "K3" "2024-12-05"
Set-TextValue "M3" ""

# Row 4
Set-TextValue "I4" "2024-11-05"
Set-TextValue "J4" ""
Set-TextValue "K4" "2024-02-05"
Set-TextValue "M4" ""
Set-TextValue "O4" ""
Set-TextValue "P4" ""

# Row 6
Set-TextValue "J6" ""
Set-TextValue "K6" "2024-11-07"
Set-TextValue "M6" ""

# Row 8
Set-TextValue "J8" ""
Set-TextValue "K8" "2024-11-07"
Set-TextValue "M8" ""
Set-TextValue "O8" ""
Set-TextValue "P8" ""

# Row 10
Set-TextValue "J10" ""
Set-TextValue "K10" "2024-11-07"
Set-TextValue "M10" ""

# Row 23
Set-TextValue "J23" ""
Set-TextValue "K23" "2024-11-11"
